$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row updates derived from the IFRS data refresh (error solve ifrs list).
# Each row entry: Set = @{ Col = Value }; Clear = @(Col, Col, ...) for cells removed entirely.
$rowEdits = @(
    @{ Row = 2; Set = @{ "D" = 209976; "E" = 11166; "F" = 11166; "G" = 11036; "H" = 8382; "I" = 8359; "J" = 23; "K" = 581836; "L" = 486448; "M" = 95388; "N" = 95254; "O" = 135; "P" = 265; "Q" = 34293; "R" = -33141; "S" = -476; "T" = 1373; "V" = 0; "W" = 5.32; "X" = 3.99; "Y" = 9.640000000000001; "Z" = 1.56; "AA" = 509.97; "AB" = 39585.88; "AC" = 16530; "AD" = 17.09; "AE" = 215654; "AF" = 1.31; "AG" = 4500; "AH" = 1.59; "AI" = 23.78; "AJ" = 47374837 }; Clear = @("U") }
    @{ Row = 3; Set = @{ "D" = 217291; "E" = 10851; "F" = 10851; "G" = 10717; "H" = 8138; "I" = 8112; "J" = 26; "K" = 632336; "L" = 531483; "M" = 100853; "N" = 100701; "O" = 151; "P" = 265; "Q" = 27012; "R" = -23575; "S" = -5203; "T" = 854; "V" = 0; "W" = 4.99; "X" = 3.74; "Y" = 8.31; "Z" = 1.34; "AA" = 526.99; "AB" = 43087.86; "AC" = 16043; "AD" = 19.17; "AE" = 234429; "AF" = 1.31; "AG" = 5150; "AI" = 27.3; "AJ" = 47374837 }; Clear = @("U") }
    @{ Row = 4; Set = @{ "D" = 216861; "E" = 10712; "F" = 10712; "G" = 11198; "H" = 8606; "I" = 8580; "J" = 27; "K" = 682175; "L" = 572093; "M" = 110082; "N" = 109908; "O" = 175; "P" = 265; "Q" = 16462; "R" = -11416; "S" = -4177; "T" = 1252; "V" = 0; "W" = 4.94; "X" = 3.97; "Y" = 8.17; "Z" = 1.31; "AA" = 519.7; "AB" = 47085.96; "AC" = 16967; "AD" = 15.82; "AE" = 258548; "AF" = 1.04; "AG" = 6100; "AH" = 2.27; "AI" = 30.23; "AJ" = 47374837 }; Clear = @("U") }
    @{ Row = 5; Set = @{ "D" = 220136; "E" = 12576; "F" = 12576; "G" = 14307; "H" = 10553; "I" = 10527; "J" = 26; "K" = 760152; "L" = 640069; "M" = 120083; "N" = 119912; "O" = 171; "P" = 265; "Q" = 25948; "R" = -20286; "S" = -2601; "T" = 769; "V" = 0; "W" = 5.71; "X" = 4.79; "Y" = 9.18; "Z" = 1.46; "AA" = 533.02; "AB" = 50863.51; "AC" = 20819; "AD" = 12.82; "AE" = 282083; "AF" = 0.95; "AG" = 10000; "AH" = 3.75; "AI" = 40.38; "AJ" = 47374837 }; Clear = @("U") }
    @{ Row = 6; Set = @{ "D" = 222090; "E" = 14543; "F" = 14543; "G" = 14599; "H" = 10733; "I" = 10705; "K" = 794164; "L" = 670333; "M" = 123830; "N" = 123640; "P" = 265; "Q" = 21262; "R" = -14527; "S" = -4256; "T" = 376; "V" = 0; "W" = 6.55; "X" = 4.83; "Y" = 8.81; "Z" = 1.38; "AA" = 541.33; "AB" = 52279.15; "AC" = 21170; "AD" = 12.71; "AE" = 290853; "AF" = 0.92; "AG" = 11500; "AH" = 4.28; "AI" = 45.67; "AJ" = 47374837 }; Clear = @("U") }
    @{ Row = 7; Set = @{ "G" = 9270; "H" = 6785; "I" = 6755; "K" = 831945; "L" = 687900; "M" = 144045; "N" = 143855; "Y" = 5.05; "Z" = 0.84; "AA" = 477.56; "AC" = 13359; "AD" = 16.84; "AE" = 338407; "AF" = 0.66; "AG" = 8103; "AH" = 3.6; "AI" = 56.83 }; Clear = @("D", "E", "P", "Q", "R", "S", "T", "U", "W", "X") }
    @{ Row = 8; Set = @{ "G" = 10260; "H" = 7445; "I" = 7435; "K" = 867455; "L" = 718515; "M" = 148940; "N" = 148750; "Y" = 5.08; "Z" = 0.88; "AA" = 482.42; "AC" = 14703; "AD" = 14.35; "AE" = 349922; "AF" = 0.6; "AG" = 9268; "AH" = 4.39; "AI" = 59.06 }; Clear = @("D", "E", "P", "Q", "R", "S", "T", "U", "W", "X") }
    @{ Row = 9; Set = @{ "G" = 11740; "H" = 8530; "I" = 8515; "K" = 899355; "L" = 744375; "M" = 154980; "N" = 154785; "Y" = 5.61; "Z" = 0.97; "AA" = 480.3; "AC" = 16839; "AD" = 12.53; "AE" = 364118; "AF" = 0.58; "AG" = 10586; "AH" = 5.02; "AI" = 58.9 }; Clear = @("D", "E", "P", "Q", "R", "S", "T", "U", "W", "X") }
)

foreach ($edit in $rowEdits) {
    $row = $edit.Row
    foreach ($col in $edit.Set.Keys) {
        $ws.Range("$col$row").Value = $edit.Set[$col]
    }
    foreach ($col in $edit.Clear) {
        $ws.Range("$col$row").ClearContents()
    }
}